$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("A5").Value = "Hiren Sojitra"
$ws.Range("B5").Value = "Amreli"
